$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# res_bus/vm_pu.xlsx case update: new vm_pu results for the 380 kV case.
# Columns B:F and I:N change for data rows 2-25 (column A index, G/H unchanged).

$bfText = @"
1.02,1.040233611731493,1.048265648958143,1.053674915278574,1.060807991849815
1.02,1.041162057858258,1.048984583468511,1.054548046946758,1.061669278766291
1.02,1.041763514732994,1.049450196544384,1.055113990706637,1.06222741147183
1.02,1.042016531159928,1.049646037867841,1.055352143864914,1.062462245675146
1.02,1.042059023309147,1.049678926163322,1.0553921442918,1.062501686769552
1.02,1.041766894908211,1.049452813004316,1.055117172017509,1.06223054857301
1.02,1.040547240784968,1.048508529471954,1.053969792561116,1.061098896480602
1.02,1.038403396128974,1.046847826228384,1.051955457942696,1.059111154768688
1.02,1.036977832637511,1.045742968486625,1.050617703549777,1.05779037939737
1.02,1.036361433022498,1.045265113317252,1.050039679393989,1.057219530390768
1.02,1.036132608006881,1.045087701711003,1.049825162251344,1.057007651911026
1.02,1.036181685700176,1.045125753226739,1.049871168429051,1.057053093285005
1.02,1.036342515550063,1.045250446673566,1.050021943520715,1.05720201317054
1.02,1.036441625821666,1.045327285743919,1.050114865827666,1.057293788913336
1.02,1.037018759668814,1.045774694002899,1.050656091211454,1.057828287147121
1.02,1.037381016617938,1.046055491256035,1.050995918669042,1.058163847533144
1.02,1.037592399675327,1.04621932903453,1.051194253068178,1.058359675950865
1.02,1.037664490166179,1.046275202506863,1.051261900072791,1.058426465583031
1.02,1.03734214110511,1.046025358823645,1.050959446099945,1.058127834545668
1.02,1.036295151479774,1.045213725199752,1.049977538845517,1.057158155538493
1.02,1.035637639078665,1.044703911975191,1.049361256153442,1.056549407902583
1.02,1.035986125302242,1.04497412633582,1.049687856225481,1.056872028010199
1.02,1.03735970701436,1.046038974207717,1.050975926117113,1.058144106950438
1.02,1.038956990648364,1.047276763380566,1.052475313646833,1.059624267974337
"@

$inText = @"
1.038683651240008,1.04532132351195,1.051025975274726,1.056420230251175,1.063533765148599,1.018995452433199
1.038839713566636,1.045895408220629,1.051557161474702,1.057106299900954,1.064209431782445,1.019188630125915
1.03893955119151,1.046266935333764,1.051900616267682,1.057550560060815,1.064646810334592,1.019313579100095
1.038981248166801,1.046423137576675,1.052044941620841,1.057737404016081,1.064830725253144,1.01936609525425
1.03898823314241,1.046449365305419,1.052069170750187,1.057768780413389,1.064861607738975,1.019374912207825
1.038940109428928,1.046269022469097,1.051902544999367,1.057553056378091,1.064649267655503,1.019314280872779
1.038736629972214,1.045515325861546,1.051205544860652,1.056652022134515,1.063762072339738,1.019060747805109
1.038369329983602,1.044187694917228,1.049975418383515,1.055066850610223,1.062200139373033,1.018613627589446
1.038118621427377,1.043302997275878,1.049154113974574,1.054011875093744,1.061159888381316,1.018315330132548
1.038008683880204,1.042920020026094,1.04879820599853,1.053555504114978,1.060709711928161,1.01818611803954
1.037967641581807,1.042777781773198,1.048665965415452,1.053386055117632,1.06054253667775,1.018138116245803
1.037976454633594,1.042808291594293,1.04869433328326,1.053422399426514,1.060578394488973,1.018148413088223
1.038005295522118,1.042908262223616,1.048787275770126,1.053541496023171,1.060695892339008,1.018182150329795
1.03802303798691,1.042969859654469,1.048844535371618,1.053614884352154,1.060768292028217,1.018202936087248
1.038125888615176,1.043328416473731,1.049177728666618,1.05404217232542,1.061189770688369,1.018323904554152
1.038190035294794,1.043553357974292,1.049386658513979,1.05431031763175,1.061454223472641,1.01839977251835
1.038227317865829,1.043684572376211,1.049508496879069,1.054466764652674,1.061608499211894,1.018444020422473
1.038240007664646,1.043729314729568,1.04955003601672,1.054520116194897,1.061661107433685,1.018459107021366
1.038183166726215,1.043529222855176,1.049364245090206,1.054281543782715,1.061425847605532,1.018391633075979
1.037996808297637,1.042878822899087,1.048759907632902,1.053506423169636,1.060661291026198,1.018172215733594
1.037878442337616,1.042469986942038,1.048379702887793,1.053019465298596,1.060180817982455,1.018034220973936
1.037941303447843,1.042686709113196,1.048581278372681,1.053277573347852,1.060435503182425,1.018107378096589
1.038186270748028,1.043540128446367,1.049374372834823,1.054294545315676,1.06143866936166,1.018395310953878
1.038465318189544,1.044530855625276,1.050293655956623,1.05547634309789,1.06260375973778,1.018729259121826
"@

$bfRows = $bfText -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }
$inRows = $inText -split "`r?`n" | Where-Object { $_.Trim().Length -gt 0 }

$nRows = $bfRows.Count
$nColsBF = ($bfRows[0] -split ",").Count
$nColsIN = ($inRows[0] -split ",").Count

$dataBF = New-Object 'object[,]' $nRows,$nColsBF
$dataIN = New-Object 'object[,]' $nRows,$nColsIN

for ($r = 0; $r -lt $nRows; $r++) {
    $bfVals = $bfRows[$r] -split ","
    for ($c = 0; $c -lt $nColsBF; $c++) {
        $dataBF[$r,$c] = $bfVals[$c] + 0
    }
    $inVals = $inRows[$r] -split ","
    for ($c = 0; $c -lt $nColsIN; $c++) {
        $dataIN[$r,$c] = $inVals[$c] + 0
    }
}

$ws.Range("B2:F25").Value = $dataBF
$ws.Range("I2:N25").Value = $dataIN
